$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{A='Entrainement'; B=45877; C='Global'; E='Romain Thunet'; F='center back'; G='01:18:46'; H=4.99; I=0.18; J=4.81; K=0.18; L=0.01; M=0; N=0; O=0; P=3.76; Q=23.71; R=4.73; S=24; T=4; U=26; V=5}
    @{A='Entrainement'; B=45877; C='Global'; E='Hedi Nasri'; F='right back'; G='01:28:54'; H=5.21; I=0.2; J=4.99; K=0.2; L=0.01; M=0; N=0; O=0; P=3.39; Q=21.67; R=5.19; S=27; T=9; U=25; V=9}
    @{A='Entrainement'; B=45877; C='Global'; E='Amine Taiar'; F='center back'; G='01:24:37'; H=5.2; I=0.18; J=5.01; K=0.17; L=0.02; M=0; N=0; O=0; P=3.58; Q=22.04; R=5.45; S=26; T=9; U=18; V=6}
    @{A='Entrainement'; B=45877; C='Global'; E='Jeremie Laurent'; F='left forward'; G='01:31:45'; H=5.79; I=0.31; J=5.47; K=0.27; L=0.05; M=0; N=0; O=0; P=3.68; Q=24.94; R=4.99; S=45; T=7; U=24; V=8}
    @{A='Entrainement'; B=45877; C='Global'; E='Ilan Ihaddadene'; F='center midfield'; G='01:29:13'; H=6.48; I=0.26; J=6.22; K=0.19; L=0.07; M=0; N=0; O=1; P=4.24; Q=25.14; R=5.28; S=30; T=8; U=22; V=5}
    @{A='Entrainement'; B=45877; C='Global'; E='Mattheo Haon'; F='right back'; G='01:22:34'; H=5.38; I=0.2; J=5.18; K=0.19; L=0.02; M=0; N=0; O=0; P=3.84; Q=22.25; R=4.97; S=28; T=8; U=24; V=8}
    @{A='Entrainement'; B=45877; C='Global'; E='Karahali Souaré'; F='right forward'; G='01:29:13'; H=5.49; I=0.27; J=5.21; K=0.24; L=0.04; M=0; N=0; O=0; P=3.45; Q=22.19; R=5.27; S=32; T=5; U=24; V=7}
    @{A='Entrainement'; B=45877; C='Global'; E='Omar Benyounes'; F='center midfield'; G='01:25:26'; H=5.52; I=0.23; J=5.29; K=0.16; L=0.06; M=0.02; N=0; O=1; P=3.71; Q=25.71; R=5.04; S=17; T=6; U=19; V=3}
    @{A='Entrainement'; B=45877; C='Global'; E='Levy Ndoutoume'; F='left back'; G='01:24:47'; H=4.57; I=0.15; J=4.41; K=0.15; L=0.01; M=0; N=0; O=0; P=3.14; Q=20.86; R=5.02; S=37; T=3; U=26; V=7}
    @{A='Entrainement'; B=45877; C='Global'; E='Emmanuel Valey'; F='left forward'; G='01:20:59'; H=5.68; I=0.22; J=5.46; K=0.18; L=0.05; M=0; N=0; O=0; P=4.14; Q=23.05; R=5.62; S=24; T=4; U=18; V=6}
)

$r = 256
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
    $r++
}

[void]$ws.Range("E269").Select()
